$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row after row 6 ("Tag model + create/query by commit") ---
# This becomes the new row 7, pushing the old rows 7-17 down to 8-18.
$ws.Rows.Item(7).EntireRow.Insert()
$ws.Range("A7").Value = "Search by commit tags "
$ws.Range("B7").Value = "To Do"

# --- 2. Insert two new rows after "Tag timeline/history view" ---
# After the previous insert, "Tag timeline/history view" (originally row 14)
# now sits at row 15. Insert two blank rows after it (at position 16), which
# become the new rows 16 and 17, pushing the Docs/Deployment/Promotion rows
# down from 15-17 to 18-20.
$ws.Rows.Item(16).EntireRow.Insert()
$ws.Rows.Item(16).EntireRow.Insert()

$ws.Range("A16").Value = "User sign on? Or is that done through the settings in the extension?"
$ws.Range("B16").Value = "To Do"

$ws.Range("A17").Value = "Modify database to support accounts!"
$ws.Range("B17").Value = "To Do"

# --- 3. Update status/value cells on existing rows (using final row numbers) ---
# "Rollback functionality" row: status -> To Do, est date -> "Needs to be testted"
$ws.Range("B4").Value = "To Do"
$ws.Range("D4").Value = "Needs to be testted"

# "Commit bug fix (creation failure)" row: status -> Complete
$ws.Range("B5").Value = "Complete"

# "Tag display on frontend" row (shifted from 7 to 8): status -> Complete
$ws.Range("B8").Value = "Complete"

# "Browser extension for ChatGPT -> ChatCommit" row (shifted from 12 to 13): status -> Working
$ws.Range("B13").Value = "Working"

# --- 4. Restore the active selection to B15 ---
$ws.Range("B15").Select()
